$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) / Volume(1h) (column E) values scraped for this run,
# keyed by row number.
$updates = @(
    @{ Row = 2; D = "68.839.86"; E = "  -0.31%  " },
    @{ Row = 3; D = "3.851.08"; E = "  +2.90%  " },
    @{ Row = 4; D = "1.00"; E = "  +0.03%  " },
    @{ Row = 5; D = "601.04"; E = "  -0.14%  " },
    @{ Row = 6; D = "162.29"; E = "  -2.85%  " },
    @{ Row = 7; D = "3.848.96"; E = "  +2.89%  " },
    @{ Row = 8; D = $null; E = "  +0.25%  " },
    @{ Row = 9; D = $null; E = "  -1.57%  " },
    @{ Row = 10; D = $null; E = "  -1.14%  " },
    @{ Row = 11; D = "6.31"; E = "  -2.28%  " },
    @{ Row = 12; D = $null; E = "  -0.22%  " },
    @{ Row = 13; D = "36.75"; E = "  -3.26%  " },
    @{ Row = 14; D = "0.0000243"; E = "  -1.94%  " },
    @{ Row = 15; D = "4.493.29"; E = "  +2.82%  " },
    @{ Row = 16; D = "3.843.45"; E = "  +2.50%  " },
    @{ Row = 17; D = "69.011.43"; E = "  -0.01%  " },
    @{ Row = 18; D = $null; E = "  +2.54%  " },
    @{ Row = 19; D = $null; E = "  +3.83%  " },
    @{ Row = 20; D = $null; E = "  -0.33%  " },
    @{ Row = 21; D = "17.13"; E = "  -1.02%  " },
    @{ Row = 22; D = "484.13"; E = "  -1.87%  " },
    @{ Row = 23; D = $null; E = "  -1.12%  " },
    @{ Row = 24; D = $null; E = "  +4.50%  " },
    @{ Row = 25; D = "83.79"; E = "  -1.20%  " },
    @{ Row = 26; D = $null; E = "  -2.36%  " },
    @{ Row = 27; D = "12.08"; E = "  -2.06%  " },
    @{ Row = 28; D = "10.00"; E = "  -0.91%  " },
    @{ Row = 29; D = $null; E = "  -0.08%  " },
    @{ Row = 30; D = $null; E = "  -0.92%  " },
    @{ Row = 31; D = "7.92"; E = "  -1.74%  " },
    @{ Row = 32; D = "4.003.14"; E = "  +2.97%  " },
    @{ Row = 33; D = $null; E = "  -3.93%  " },
    @{ Row = 34; D = "32.16"; E = "  +2.01%  " },
    @{ Row = 35; D = "3.799.59"; E = "  +3.26%  " },
    @{ Row = 36; D = $null; E = "  -1.79%  " },
    @{ Row = 37; D = $null; E = "  +1.52%  " },
    @{ Row = 38; D = "0.139"; E = "  +3.50%  " },
    @{ Row = 39; D = "5.88"; E = "  -1.22%  " },
    @{ Row = 40; D = "1.00"; E = "  +0.09%  " },
    @{ Row = 41; D = $null; E = "  -1.99%  " },
    @{ Row = 42; D = "438.11"; E = "  +1.52%  " },
    @{ Row = 43; D = "2.97"; E = "  -0.68%  " },
    @{ Row = 44; D = $null; E = "  -0.55%  " },
    @{ Row = 45; D = $null; E = "  -1.03%  " },
    @{ Row = 46; D = $null; E = "  +0.00%  " },
    @{ Row = 47; D = "8.38"; E = "  -1.15%  " },
    @{ Row = 48; D = "26.44"; E = "  +11.55%  " },
    @{ Row = 49; D = "142.86"; E = "  +1.12%  " },
    @{ Row = 50; D = "2.830.85"; E = "  +1.64%  " },
    @{ Row = 51; D = "0.0357"; E = "  +1.44%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D$($u.Row)")
        # Force plain-text storage so price strings such as "1.00" or
        # "10.00" are not silently coerced into numbers (losing trailing
        # zeros / the literal text representation scraped from the site).
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
